$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Dyson Daniels, PG,SG, Atlanta Hawks -> Kevin Huerter, SG,SF, Sacramento Kings
$ws.Range("A3").Value = "Kevin Huerter"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Sacramento Kings"

# Row 7: Naz Reid, PF,C, Minnesota Timberwolves -> Kyle Filipowski, PF,C, Utah Jazz
$ws.Range("A7").Value = "Kyle Filipowski"
$ws.Range("C7").Value = "Utah Jazz"

# Row 14: Naji Marshall, SG,SF, Dallas Mavericks -> Dyson Daniels, PG,SG, Atlanta Hawks
$ws.Range("A14").Value = "Dyson Daniels"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Atlanta Hawks"

# Row 15: Kyle Filipowski, PF, Utah Jazz -> Naz Reid, PF,C, Minnesota Timberwolves
$ws.Range("A15").Value = "Naz Reid"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Minnesota Timberwolves"
